$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 4,65
$data[0,0] = "2025-09-01 09:25:46"
$data[0,1] = "data\Multi-Illumination"
$data[0,2] = "a-1-b-1"
$data[0,3] = 1848754
$data[0,4] = 0.002
$data[0,5] = 0.004
$data[0,6] = 7
$data[0,7] = 0.000003786333930852888
$data[0,8] = 0.9999962136660692
$data[0,9] = 1848747
$data[0,10] = 58.06580200000001
$data[0,11] = 0.201039354048
$data[0,12] = 1831404
$data[0,13] = 64.48812999999998
$data[0,14] = 0.16406526557
$data[0,15] = -0.004257
$data[0,16] = 0.004802
$data[0,17] = 0.00003140819268401788
$data[0,18] = 0.000044
$data[0,19] = 0.0003297629039704233
$data[0,20] = 0.0003282637632565164
$data[0,21] = 0.0002423816000783233
$data[0,22] = 0.0002683506
$data[0,23] = -0.0009890000000000001
$data[0,24] = 0.0009890000000000001
$data[0,25] = 0.00003521239988555228
$data[0,26] = 0.000045
$data[0,27] = 0.000299306565652362
$data[0,28] = 0.0002972280389480636
$data[0,29] = 0.0002315634911794449
$data[0,30] = 0.0002653854
$data[0,31] = 3
$data[0,32] = 0.00098928871191127
$data[0,33] = "rmse"
$data[0,34] = 1831404
$data[0,35] = 1041360
$data[0,36] = 787214
$data[0,37] = 6541
$data[0,38] = 10802
$data[0,39] = 17343
$data[0,40] = -0.0003703124026984951
$data[0,41] = 0.001412373977284256
$data[0,42] = -0.000513
$data[0,43] = -0.000141
$data[0,44] = 0.000222
$data[0,45] = 0.000523
$data[0,46] = 0.000363
$data[0,47] = -0.000491
$data[0,48] = -0.000137
$data[0,49] = 0.000221
$data[0,50] = 0.000513
$data[0,51] = 0.000358
$data[0,52] = 0.00003140819268401788
$data[0,53] = 0.0003282637632565164
$data[0,54] = 1127175200681164
$data[0,55] = 12.22131184167564
$data[0,56] = 0.00443167676357952
$data[0,57] = -0.004257395391128698
$data[0,58] = 0.0001434337456783985
$data[0,59] = -0.7175292079112111
$data[0,60] = 14923656066230.08
$data[0,61] = -0.5887868609417466
$data[0,62] = 5.582690835825172
$data[0,63] = "data\Multi-Illumination\python_a-1-b-1_m3c2_distances.txt"
$data[0,64] = "data\Multi-Illumination\python_a-1-b-1_m3c2_params.txt"
$data[1,0] = "2025-09-01 09:27:16"
$data[1,1] = "data\Multi-Illumination"
$data[1,2] = "a-1-b-1-AI"
$data[1,3] = 1848754
$data[1,4] = 0.002
$data[1,5] = 0.004
$data[1,6] = 3353
$data[1,7] = 0.001813653952878533
$data[1,8] = 0.9981863460471214
$data[1,9] = 1845401
$data[1,10] = -671.42272
$data[1,11] = 1.611225486426
$data[1,12] = 1825744
$data[1,13] = -626.261582
$data[1,14] = 1.407286543668
$data[1,15] = -0.00432
$data[1,16] = 0.004148
$data[1,17] = -0.0003638356758233035
$data[1,18] = -0.000283
$data[1,19] = 0.0009343999051170693
$data[1,20] = 0.0008606548574666783
$data[1,21] = 0.0006962032219555533
$data[1,22] = 0.000748713
$data[1,23] = -0.002803
$data[1,24] = 0.002803
$data[1,25] = -0.0003430171929909122
$data[1,26] = -0.000275
$data[1,27] = 0.0008779531198573015
$data[1,28] = 0.0008081713221711127
$data[1,29] = 0.0006691691233820295
$data[1,30] = 0.0007383347999999999
$data[1,31] = 3
$data[1,32] = 0.002803199715351208
$data[1,33] = "rmse"
$data[1,34] = 1825744
$data[1,35] = 624201
$data[1,36] = 1200506
$data[1,37] = 2784
$data[1,38] = 16873
$data[1,39] = 19657
$data[1,40] = -0.002297458310016787
$data[1,41] = 0.002257556550696858
$data[1,42] = -0.001857
$data[1,43] = -0.000859
$data[1,44] = 0.00017
$data[1,45] = 0.000891
$data[1,46] = 0.001029
$data[1,47] = -0.00177
$data[1,48] = -0.00084
$data[1,49] = 0.000172
$data[1,50] = 0.000882
$data[1,51] = 0.001012
$data[1,52] = -0.0003638356758233035
$data[1,53] = 0.0008606548574666783
$data[1,54] = 248623.2555508183
$data[1,55] = 5.100174121428031
$data[1,56] = 0.004338516079674972
$data[1,57] = -0.004368371937607233
$data[1,58] = -0.0002115922981897056
$data[1,59] = -0.2677860148641807
$data[1,60] = 1326327888.537201
$data[1,61] = -0.2770356403359415
$data[1,62] = 1.312988890863815
$data[1,63] = "data\Multi-Illumination\python_a-1-b-1-AI_m3c2_distances.txt"
$data[1,64] = "data\Multi-Illumination\python_a-1-b-1-AI_m3c2_params.txt"
$data[2,0] = "2025-09-01 09:29:02"
$data[2,1] = "data\Multi-Illumination"
$data[2,2] = "a-1-AI-b-1"
$data[2,3] = 2143284
$data[2,4] = 0.002
$data[2,5] = 0.004
$data[2,6] = 214008
$data[2,7] = 0.09985050977845213
$data[2,8] = 0.9001494902215479
$data[2,9] = 1929276
$data[2,10] = 944.0018060000002
$data[2,11] = 2.42405412044
$data[2,12] = 1917090
$data[2,13] = 912.6991869999999
$data[2,14] = 2.261979214093
$data[2,15] = -0.005808
$data[2,16] = 0.004936
$data[2,17] = 0.0004893036589891754
$data[2,18] = 0.000482
$data[2,19] = 0.001120918340191286
$data[2,20] = 0.001008483938730307
$data[2,21] = 0.0008716386053628407
$data[2,22] = 0.0008776991999999999
$data[2,23] = -0.003362
$data[2,24] = 0.003362
$data[2,25] = 0.0004760857273263122
$data[2,26] = 0.000477
$data[2,27] = 0.001086233153461209
$data[2,28] = 0.000976342585322621
$data[2,29] = 0.000854039257937812
$data[2,30] = 0.0008702861999999999
$data[2,31] = 3
$data[2,32] = 0.003362755020573859
$data[2,33] = "rmse"
$data[2,34] = 1917090
$data[2,35] = 1367184
$data[2,36] = 549002
$data[2,37] = 10429
$data[2,38] = 1757
$data[2,39] = 12186
$data[2,40] = 0.002568736172657148
$data[2,41] = 0.002588761337714701
$data[2,42] = -0.001169
$data[2,43] = -0.000098
$data[2,44] = 0.001086
$data[2,45] = 0.002158
$data[2,46] = 0.001184
$data[2,47] = -0.001162
$data[2,48] = -0.0001
$data[2,49] = 0.001074
$data[2,50] = 0.002095
$data[2,51] = 0.001174
$data[2,52] = 0.0004893036589891754
$data[2,53] = 0.001008483938730307
$data[2,54] = 122338.8498843068
$data[2,55] = 6.632135228644302
$data[2,56] = 0.006722693953013577
$data[2,57] = -0.005808409162356089
$data[2,58] = 0.0007506396059552209
$data[2,59] = -0.4327674980723418
$data[2,60] = 203600348.2673969
$data[2,61] = -0.04436061331310492
$data[2,62] = 0.9798550044227512
$data[2,63] = "data\Multi-Illumination\python_a-1-AI-b-1_m3c2_distances.txt"
$data[2,64] = "data\Multi-Illumination\python_a-1-AI-b-1_m3c2_params.txt"
$data[3,0] = "2025-09-01 09:31:04"
$data[3,1] = "data\Multi-Illumination"
$data[3,2] = "a-1-AI-b-1-AI"
$data[3,3] = 2143284
$data[3,4] = 0.002
$data[3,5] = 0.004
$data[3,6] = 25924
$data[3,7] = 0.01209545725158215
$data[3,8] = 0.9879045427484179
$data[3,9] = 2117360
$data[3,10] = 522.5222369999999
$data[3,11] = 1.987969683329
$data[3,12] = 2067166
$data[3,13] = 387.117282
$data[3,14] = 1.390636687952
$data[3,15] = -0.005688
$data[3,16] = 0.00617
$data[3,17] = 0.0002467800643253863
$data[3,18] = 0.000183
$data[3,19] = 0.000968963742060044
$data[3,20] = 0.0009370113837507853
$data[3,21] = 0.0006621369068084785
$data[3,22] = 0.0006063833999999999
$data[3,23] = -0.002906
$data[3,24] = 0.002906
$data[3,25] = 0.0001872695671271683
$data[3,26] = 0.00017
$data[3,27] = 0.0008201988665263999
$data[3,28] = 0.0007985338376544816
$data[3,29] = 0.0005949083663334245
$data[3,30] = 0.0005841444
$data[3,31] = 3
$data[3,32] = 0.002906891226180132
$data[3,33] = "rmse"
$data[3,34] = 2067166
$data[3,35] = 1281193
$data[3,36] = 784476
$data[3,37] = 44868
$data[3,38] = 5326
$data[3,39] = 50194
$data[3,40] = 0.002697632286727497
$data[3,41] = 0.002150178147824242
$data[3,42] = -0.001117
$data[3,43] = -0.000234
$data[3,44] = 0.0005820000000000001
$data[3,45] = 0.002048
$data[3,46] = 0.0008160000000000001
$data[3,47] = -0.0011
$data[3,48] = -0.000241
$data[3,49] = 0.000545
$data[3,50] = 0.001659
$data[3,51] = 0.000786
$data[3,52] = 0.0002467800643253863
$data[3,53] = 0.0009370113837507853
$data[3,54] = 1561851.362796817
$data[3,55] = 5.939962160905105
$data[3,56] = 0.006340247401192347
$data[3,57] = -0.005688138567817923
$data[3,58] = 0.0004583629744138711
$data[3,59] = -0.3670612797172161
$data[3,60] = 11693423546.35915
$data[3,61] = 0.7472737316746452
$data[3,62] = 2.994743908023999
$data[3,63] = "data\Multi-Illumination\python_a-1-AI-b-1-AI_m3c2_distances.txt"
$data[3,64] = "data\Multi-Illumination\python_a-1-AI-b-1-AI_m3c2_params.txt"

$ws.Range("A28:BM31").Value = $data
